$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# === Step 1: introduce new shared strings in original authoring order ===
$ws.Range("A376").Value = "건강"
$ws.Range("B376").Value = "здоровье"
$ws.Range("D376").Value = "condition"
$ws.Range("E376").Value = "состояние"
$ws.Range("A377").Value = "머리"
$ws.Range("A378").Value = "귀"
$ws.Range("A379").Value = "다리"
$ws.Range("A380").Value = "얼굴"
$ws.Range("A382").Value = "발"
$ws.Range("A388").Value = "허리"
$ws.Range("A389").Value = "입"
$ws.Range("A390").Value = "손"
$ws.Range("A391").Value = "어깨"
$ws.Range("B377").Value = "голово (волосы)"
$ws.Range("B381").Value = "зуб"
$ws.Range("B383").Value = "глаз"
$ws.Range("B384").Value = "шея"
$ws.Range("B385").Value = "живот"
$ws.Range("B386").Value = "нос"
$ws.Range("B387").Value = "рука"
$ws.Range("B389").Value = "рот"
$ws.Range("B391").Value = "плечи"
$ws.Range("B380").Value = "лицо"
$ws.Range("B382").Value = "стопа"
$ws.Range("F377").Value = "body_parts"
$ws.Range("G377").Value = "части тела"
$ws.Range("A392").Value = "참다"
$ws.Range("B392").Value = "терпеть"
$ws.Range("A393").Value = "건강하다"
$ws.Range("A395").Value = "감기에 걸리다"
$ws.Range("A396").Value = "닛다"
$ws.Range("A397").Value = "몸이 안 졸다"
$ws.Range("A398").Value = "열이나다"
$ws.Range("A399").Value = "좋아지다"
$ws.Range("A400").Value = "다치다"
$ws.Range("A401").Value = "기침을 하다"
$ws.Range("B393").Value = "здоровый"
$ws.Range("B396").Value = "становится лучше"
$ws.Range("B397").Value = "плохо себя чувствую"
$ws.Range("B398").Value = "температура"
$ws.Range("B400").Value = "пораниться"
$ws.Range("B401").Value = "кашлять"
$ws.Range("B379").Value = "ноги"
$ws.Range("B378").Value = "уши"
$ws.Range("B388").Value = "спина (поясница)"
$ws.Range("B390").Value = "кисти"
$ws.Range("B399").Value = "улучшаться (попровляться)"
$ws.Range("B395").Value = "простудиться"
$ws.Range("B394").Value = "болеть"
$ws.Range("B402").Value = "секрет"
$ws.Range("B403").Value = "рано"
$ws.Range("B404").Value = "вдоволь (высыпаться)"
$ws.Range("B405").Value = "до, перед"
$ws.Range("B406").Value = "менять"
$ws.Range("B407").Value = "слово"
$ws.Range("B408").Value = "быстро"
$ws.Range("B409").Value = "лекарство"
$ws.Range("B410").Value = "мыть / мыться"
$ws.Range("B411").Value = "прививка"
$ws.Range("A402").Value = "비밀"
$ws.Range("A403").Value = "일찍"
$ws.Range("A404").Value = "푹"
$ws.Range("A405").Value = "전"
$ws.Range("A406").Value = "바꾸다"
$ws.Range("A407").Value = "말"
$ws.Range("A408").Value = "빨리"
$ws.Range("A409").Value = "약"
$ws.Range("A410").Value = "씻다"
$ws.Range("A411").Value = "습관"
$ws.Range("D231").Value = "annother"

# === Step 2: fill remaining cells for modified existing rows (reuse existing strings) ===
$ws.Range("E231").Value = "другое"
$ws.Range("D232").Value = "annother"
$ws.Range("E232").Value = "другое"
$ws.Range("D236").Value = "people"
$ws.Range("E236").Value = "люди"
$ws.Range("D245").Value = "annother"
$ws.Range("E245").Value = "другое"
$ws.Range("D253").Value = "annother"
$ws.Range("E253").Value = "другое"
$ws.Range("D255").Value = "annother"
$ws.Range("E255").Value = "другое"

# === Step 3: fill remaining cells for new rows 376-411 (numbers + already-written text cols) ===
$ws.Range("J376").Value = "noun"
$ws.Range("K376").Value = 2
$ws.Range("L376").Value = 10
$ws.Range("M376").Value = "words"
$ws.Range("D377").Value = "people"
$ws.Range("E377").Value = "люди"
$ws.Range("J377").Value = "noun"
$ws.Range("K377").Value = 2
$ws.Range("L377").Value = 10
$ws.Range("M377").Value = "words"
$ws.Range("D378").Value = "people"
$ws.Range("E378").Value = "люди"
$ws.Range("F378").Value = "body_parts"
$ws.Range("G378").Value = "части тела"
$ws.Range("J378").Value = "noun"
$ws.Range("K378").Value = 2
$ws.Range("L378").Value = 10
$ws.Range("M378").Value = "words"
$ws.Range("D379").Value = "people"
$ws.Range("E379").Value = "люди"
$ws.Range("F379").Value = "body_parts"
$ws.Range("G379").Value = "части тела"
$ws.Range("J379").Value = "noun"
$ws.Range("K379").Value = 2
$ws.Range("L379").Value = 10
$ws.Range("M379").Value = "words"
$ws.Range("D380").Value = "people"
$ws.Range("E380").Value = "люди"
$ws.Range("F380").Value = "body_parts"
$ws.Range("G380").Value = "части тела"
$ws.Range("J380").Value = "noun"
$ws.Range("K380").Value = 2
$ws.Range("L380").Value = 10
$ws.Range("M380").Value = "words"
$ws.Range("A381").Value = "이"
$ws.Range("D381").Value = "people"
$ws.Range("E381").Value = "люди"
$ws.Range("F381").Value = "body_parts"
$ws.Range("G381").Value = "части тела"
$ws.Range("J381").Value = "noun"
$ws.Range("K381").Value = 2
$ws.Range("L381").Value = 10
$ws.Range("M381").Value = "words"
$ws.Range("D382").Value = "people"
$ws.Range("E382").Value = "люди"
$ws.Range("F382").Value = "body_parts"
$ws.Range("G382").Value = "части тела"
$ws.Range("J382").Value = "noun"
$ws.Range("K382").Value = 2
$ws.Range("L382").Value = 10
$ws.Range("M382").Value = "words"
$ws.Range("A383").Value = "눈"
$ws.Range("D383").Value = "people"
$ws.Range("E383").Value = "люди"
$ws.Range("F383").Value = "body_parts"
$ws.Range("G383").Value = "части тела"
$ws.Range("J383").Value = "noun"
$ws.Range("K383").Value = 2
$ws.Range("L383").Value = 10
$ws.Range("M383").Value = "words"
$ws.Range("A384").Value = "목"
$ws.Range("D384").Value = "people"
$ws.Range("E384").Value = "люди"
$ws.Range("F384").Value = "body_parts"
$ws.Range("G384").Value = "части тела"
$ws.Range("J384").Value = "noun"
$ws.Range("K384").Value = 2
$ws.Range("L384").Value = 10
$ws.Range("M384").Value = "words"
$ws.Range("A385").Value = "배"
$ws.Range("D385").Value = "people"
$ws.Range("E385").Value = "люди"
$ws.Range("F385").Value = "body_parts"
$ws.Range("G385").Value = "части тела"
$ws.Range("J385").Value = "noun"
$ws.Range("K385").Value = 2
$ws.Range("L385").Value = 10
$ws.Range("M385").Value = "words"
$ws.Range("A386").Value = "코"
$ws.Range("D386").Value = "people"
$ws.Range("E386").Value = "люди"
$ws.Range("F386").Value = "body_parts"
$ws.Range("G386").Value = "части тела"
$ws.Range("J386").Value = "noun"
$ws.Range("K386").Value = 2
$ws.Range("L386").Value = 10
$ws.Range("M386").Value = "words"
$ws.Range("A387").Value = "팔다"
$ws.Range("D387").Value = "people"
$ws.Range("E387").Value = "люди"
$ws.Range("F387").Value = "body_parts"
$ws.Range("G387").Value = "части тела"
$ws.Range("J387").Value = "noun"
$ws.Range("K387").Value = 2
$ws.Range("L387").Value = 10
$ws.Range("M387").Value = "words"
$ws.Range("D388").Value = "people"
$ws.Range("E388").Value = "люди"
$ws.Range("F388").Value = "body_parts"
$ws.Range("G388").Value = "части тела"
$ws.Range("J388").Value = "noun"
$ws.Range("K388").Value = 2
$ws.Range("L388").Value = 10
$ws.Range("M388").Value = "words"
$ws.Range("D389").Value = "people"
$ws.Range("E389").Value = "люди"
$ws.Range("F389").Value = "body_parts"
$ws.Range("G389").Value = "части тела"
$ws.Range("J389").Value = "noun"
$ws.Range("K389").Value = 2
$ws.Range("L389").Value = 10
$ws.Range("M389").Value = "words"
$ws.Range("D390").Value = "people"
$ws.Range("E390").Value = "люди"
$ws.Range("F390").Value = "body_parts"
$ws.Range("G390").Value = "части тела"
$ws.Range("J390").Value = "noun"
$ws.Range("K390").Value = 2
$ws.Range("L390").Value = 10
$ws.Range("M390").Value = "words"
$ws.Range("D391").Value = "people"
$ws.Range("E391").Value = "люди"
$ws.Range("F391").Value = "body_parts"
$ws.Range("G391").Value = "части тела"
$ws.Range("J391").Value = "noun"
$ws.Range("K391").Value = 2
$ws.Range("L391").Value = 10
$ws.Range("M391").Value = "words"
$ws.Range("D392").Value = "condition"
$ws.Range("E392").Value = "состояние"
$ws.Range("J392").Value = "verb"
$ws.Range("K392").Value = 2
$ws.Range("L392").Value = 10
$ws.Range("M392").Value = "expression"
$ws.Range("D393").Value = "condition"
$ws.Range("E393").Value = "состояние"
$ws.Range("J393").Value = "verb"
$ws.Range("K393").Value = 2
$ws.Range("L393").Value = 10
$ws.Range("M393").Value = "words"
$ws.Range("A394").Value = "아프다"
$ws.Range("D394").Value = "condition"
$ws.Range("E394").Value = "состояние"
$ws.Range("J394").Value = "noun"
$ws.Range("K394").Value = 2
$ws.Range("L394").Value = 10
$ws.Range("M394").Value = "words"
$ws.Range("D395").Value = "condition"
$ws.Range("E395").Value = "состояние"
$ws.Range("J395").Value = "verb"
$ws.Range("K395").Value = 2
$ws.Range("L395").Value = 10
$ws.Range("M395").Value = "expression"
$ws.Range("D396").Value = "condition"
$ws.Range("E396").Value = "состояние"
$ws.Range("K396").Value = 2
$ws.Range("L396").Value = 10
$ws.Range("M396").Value = "expression"
$ws.Range("D397").Value = "condition"
$ws.Range("E397").Value = "состояние"
$ws.Range("K397").Value = 2
$ws.Range("L397").Value = 10
$ws.Range("M397").Value = "expression"
$ws.Range("D398").Value = "condition"
$ws.Range("E398").Value = "состояние"
$ws.Range("J398").Value = "noun"
$ws.Range("K398").Value = 2
$ws.Range("L398").Value = 10
$ws.Range("M398").Value = "expression"
$ws.Range("D399").Value = "condition"
$ws.Range("E399").Value = "состояние"
$ws.Range("J399").Value = "verb"
$ws.Range("K399").Value = 2
$ws.Range("L399").Value = 10
$ws.Range("M399").Value = "expression"
$ws.Range("D400").Value = "condition"
$ws.Range("E400").Value = "состояние"
$ws.Range("J400").Value = "verb"
$ws.Range("K400").Value = 2
$ws.Range("L400").Value = 10
$ws.Range("M400").Value = "expression"
$ws.Range("D401").Value = "condition"
$ws.Range("E401").Value = "состояние"
$ws.Range("J401").Value = "verb"
$ws.Range("K401").Value = 2
$ws.Range("L401").Value = 10
$ws.Range("M401").Value = "words"
$ws.Range("D402").Value = "annother"
$ws.Range("E402").Value = "другое"
$ws.Range("J402").Value = "noun"
$ws.Range("K402").Value = 2
$ws.Range("L402").Value = 10
$ws.Range("M402").Value = "words"
$ws.Range("D403").Value = "time"
$ws.Range("E403").Value = "время"
$ws.Range("J403").Value = "adverb"
$ws.Range("K403").Value = 2
$ws.Range("L403").Value = 10
$ws.Range("M403").Value = "words"
$ws.Range("D404").Value = "annother"
$ws.Range("E404").Value = "другое"
$ws.Range("J404").Value = "adverb"
$ws.Range("K404").Value = 2
$ws.Range("L404").Value = 10
$ws.Range("M404").Value = "words"
$ws.Range("D405").Value = "annother"
$ws.Range("E405").Value = "другое"
$ws.Range("J405").Value = "adverb"
$ws.Range("K405").Value = 2
$ws.Range("L405").Value = 10
$ws.Range("M405").Value = "words"
$ws.Range("D406").Value = "action"
$ws.Range("E406").Value = "действие"
$ws.Range("J406").Value = "verb"
$ws.Range("K406").Value = 2
$ws.Range("L406").Value = 10
$ws.Range("M406").Value = "words"
$ws.Range("D407").Value = "annother"
$ws.Range("E407").Value = "другое"
$ws.Range("J407").Value = "noun"
$ws.Range("K407").Value = 2
$ws.Range("L407").Value = 10
$ws.Range("M407").Value = "words"
$ws.Range("D408").Value = "annother"
$ws.Range("E408").Value = "другое"
$ws.Range("J408").Value = "adverb"
$ws.Range("K408").Value = 2
$ws.Range("L408").Value = 10
$ws.Range("M408").Value = "words"
$ws.Range("D409").Value = "condition"
$ws.Range("E409").Value = "состояние"
$ws.Range("J409").Value = "noun"
$ws.Range("K409").Value = 2
$ws.Range("L409").Value = 10
$ws.Range("M409").Value = "words"
$ws.Range("D410").Value = "action"
$ws.Range("E410").Value = "действие"
$ws.Range("J410").Value = "verb"
$ws.Range("K410").Value = 2
$ws.Range("L410").Value = 10
$ws.Range("M410").Value = "words"
$ws.Range("D411").Value = "condition"
$ws.Range("E411").Value = "состояние"
$ws.Range("J411").Value = "noun"
$ws.Range("K411").Value = 2
$ws.Range("L411").Value = 10
$ws.Range("M411").Value = "words"

# === Step 4: worksheet-level formatting / view changes ===
$ws.Columns.Item(7).ColumnWidth = 15.2

# Reset and reapply AutoFilter over the new extended range
$ws.AutoFilterMode = $false
$ws.Range("A1:M411").AutoFilter()

# Update the _FilterDatabase defined name to match the new range
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Лист1!`$A`$1:`$M`$411"

# Zoom + selection/pane scroll position
$excel.ActiveWindow.Zoom = 115
$ws.Range("E5").Select()
